# Update "Pagos" (column F) and "Inscrições homologadas" (column H) values
# for the rows whose data changed, per the commit diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

# Map of row number -> new F (Pagos) and H (Inscrições homologadas) values
$updates = @{
    2  = @{ F = 27;  H = 29 }
    3  = @{ F = 24;  H = 27 }
    4  = @{ F = 18;  H = 21 }
    5  = @{ F = 25;  H = 29 }
    6  = @{ F = 36;  H = 43 }
    8  = @{ F = 22;  H = 28 }
    9  = @{ F = 10;  H = 14 }
    10 = @{ F = 22;  H = 24 }
    11 = @{ F = 19;  H = 20 }
    12 = @{ F = 31;  H = 33 }
    14 = @{ F = 27;  H = 29 }
    15 = @{ F = 85;  H = 96 }
    16 = @{ F = 133; H = 221 }
    17 = @{ F = 23;  H = 24 }
    18 = @{ F = 48;  H = 71 }
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row].F
    $ws.Range("H$row").Value = $updates[$row].H
}
